$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 99.99 -> 0M
$t.Cell(1, 1).Range.Text = "0M"

# Row 2: 0.45 -> 0M
$t.Cell(2, 1).Range.Text = "0M"

# Row 3: 3980 -> 0M
$t.Cell(3, 1).Range.Text = "0M"

# Row 4: 609 -> 2436
$t.Cell(4, 1).Range.Text = "2436"

# Row 6: 0.00033 -> 0.00199
$t.Cell(6, 1).Range.Text = "0.00199"

# Row 7: 0.00009 -> 0.00019
$t.Cell(7, 1).Range.Text = "0.00019"

# Row 8: 0.00002 -> 0.00007
$t.Cell(8, 1).Range.Text = "0.00007"

# Row 9: 0.00008 -> 0.00027
$t.Cell(9, 1).Range.Text = "0.00027"

# Row 10: 0.00009 -> 0.00030
$t.Cell(10, 1).Range.Text = "0.00030"

# Row 11: 0.00009 -> 0.00039
$t.Cell(11, 1).Range.Text = "0.00039"

# Row 12: 0.05401 -> 0.45349
$t.Cell(12, 1).Range.Text = "0.45349"

# Row 44: tab-separated list collapses to 99.99
$t.Cell(44, 1).Range.Text = "99.99"

# Row 45: tab-separated list collapses to 0.45
$t.Cell(45, 1).Range.Text = "0.45"

# Row 46: tab-separated list collapses to 3980
$t.Cell(46, 1).Range.Text = "3980"
